$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("var clinker kiln")

$ws.Range("F3").Value = 0.8
$ws.Range("F4").Value = 0.8
$ws.Range("F5").Value = 0.8
$ws.Range("F6").Value = 0.8

$ws.Activate()
$ws.Range("F7").Select()
